$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3992230000000001
$ws.Range("H2").Value = 1.197669
$ws.Range("N2").Value = 0.420294
$ws.Range("O2").Value = 0.4537347012141868
$ws.Range("P2").Value = 0.4537347012141868
$ws.Range("Q2").Value = 0.05593034385400001
$ws.Range("R2").Value = 0.503373094686
$ws.Range("S2").Value = 0.4537347012141868
$ws.Range("T2").Value = 0.4537347012141868

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3992230000000001
$ws.Range("H3").Value = 1.197669
$ws.Range("O3").Value = 0.270455867921697
$ws.Range("P3").Value = 0.270455867921697
$ws.Range("Q3").Value = 0.03333818120966667
$ws.Range("R3").Value = 0.300043630887
$ws.Range("S3").Value = 0.270455867921697
$ws.Range("T3").Value = 0.270455867921697

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3992230000000001
$ws.Range("H4").Value = 1.197669
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08516066666666666
$ws.Range("N4").Value = 0.255482
$ws.Range("O4").Value = 0.2758094308641162
$ws.Range("P4").Value = 0.2758094308641162
$ws.Range("Q4").Value = 0.03399809682866667
$ws.Range("R4").Value = 0.305982871458
$ws.Range("S4").Value = 0.2758094308641162
$ws.Range("T4").Value = 0.2758094308641162
